$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header for column E ("Test Number") and copy D1's formatting (bold/centered header style)
$ws.Range("E1").Value = "Test Number"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "Test Number" values for the last two rows
$ws.Range("E16").Value = 5035
$ws.Range("E17").Value = 7286

# Widen columns B:E to fit the new column (closest achievable width)
$ws.Range("B1:E1").ColumnWidth = 13.571428571428571

# Move the active selection to E3, matching the post-edit workbook state
$ws.Range("E3").Select() | Out-Null
